# Updated test data as per new implementation:
#  - Loading Details Name "Main Processor 24V (A)" (L8) is renamed to "24V Rail(A)"
#  - A review comment from Alpesh Dhakad is added to H8 ("Printer 2") noting that
#    Printer 2 is not visible and PLX800-E shows instead
#  - Selection moves to H8 (the commented cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Main Processor 24V (A)" loading-details column header to "24V Rail(A)"
$ws.Range("L8").Value = "24V Rail(A)"

# Leave a review comment on H8 ("Printer 2") documenting the observed defect
$commentText = "Alpesh Dhakad:" + [char]10 + "Printer 2 is not visible. Instead of it PLX800-E is visible." + [char]10
$ws.Range("H8").AddComment($commentText)

# Reflect that H8 is now the cell of interest
$ws.Range("H8").Select()
